# Adds "2022-Q4" quarterly data:
#  1. Inserts a new row into the "总计" (summary) sheet for the 2022-Q4 totals,
#     shifting the existing 2022-Q3 / 2022-Q2 / 2022-Q1 rows down by one.
#  2. Inserts a brand-new "2022-Q4" worksheet (positioned right after "总计",
#     before "2022-Q3") and fills it with the fund holdings detail for the
#     quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet - add the 2022-Q4 summary row just below the header row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Match the bordered/centered style already used by the other A-column cells.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 23
$summary.Cells.Item(2,4).Value = 0.45

# Column A is a running 0-based row index (not a fixed label), so the rows
# that shifted down need their index bumped by one to stay sequential.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3

# ---------------------------------------------------------------------------
# 2. Brand new "2022-Q4" worksheet with the fund-level detail, inserted
#    immediately before the existing "2022-Q3" tab.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

$rows = @(
    @(0,  "012846", "恒越蓝筹精选混合",               "5.90", "86.65", "3.58", "0.2112", 4),
    @(1,  "009379", "中银成长优选股票A",               "1.51", "83.26", "3.71", "0.0560", 5),
    @(2,  "003131", "国寿安保强国智造灵活配置混合",     "1.38", "92.54", "1.75", "0.0242", 8),
    @(3,  "004332", "恒生前海沪港深新兴产业精选混合",   "0.49", "75.81", "4.23", "0.0207", 5),
    @(4,  "011735", "国寿安保裕丰混合C",               "2.64", "27.90", "0.66", "0.0174", 8),
    @(5,  "010206", "国寿安保裕安混合C",               "2.00", "33.78", "0.71", "0.0142", 10),
    @(6,  "012461", "西藏东财国证龙头家电指数A",       "0.63", "94.76", "2.14", "0.0135", 9),
    @(7,  "010205", "国寿安保裕安混合A",               "1.89", "33.78", "0.71", "0.0134", 10),
    @(8,  "012462", "西藏东财国证龙头家电指数C",       "0.56", "94.76", "2.14", "0.0120", 9),
    @(9,  "013054", "天弘国证龙头家电指数C",           "0.41", "95.30", "2.16", "0.0089", 9),
    @(10, "011734", "国寿安保裕丰混合A",               "1.22", "27.90", "0.66", "0.0081", 8),
    @(11, "159730", "博时国证龙头家电ETF",             "0.35", "99.13", "2.25", "0.0079", 9),
    @(12, "010765", "国寿安保华丰混合A",               "0.40", "83.84", "1.86", "0.0074", 10),
    @(13, "013383", "恒生前海高端制造混合A",           "0.11", "84.98", "5.82", "0.0064", 3),
    @(14, "080007", "长盛同鑫行业配置混合A",           "0.20", "84.46", "3.18", "0.0064", 2),
    @(15, "014455", "中银成长优选股票C",               "0.12", "83.26", "3.71", "0.0045", 5),
    @(16, "010487", "中银顺盈回报一年持有期混合",       "0.75", "21.31", "0.58", "0.0044", 8),
    @(17, "080015", "长盛中小盘精选混合",               "0.13", "84.17", "2.94", "0.0038", 3),
    @(18, "013053", "天弘国证龙头家电指数A",           "0.17", "95.30", "2.16", "0.0037", 9),
    @(19, "013384", "恒生前海高端制造混合C",           "0.04", "84.98", "5.82", "0.0023", 3),
    @(20, "001932", "国寿安保灵活优选混合",             "0.11", "39.50", "1.10", "0.0012", 5),
    @(21, "010991", "长盛同鑫行业配置混合C",           "0.02", "84.46", "3.18", "0.0006", 2),
    @(22, "010766", "国寿安保华丰混合C",               "0.01", "83.84", "1.86", "0.0002", 10)
)

# Columns B (fund code) and D/E/F/G (percentages & scale figures) are stored
# as text in the source data (to keep leading zeros / fixed decimal places),
# so force a "Text" number format on those specific cells before writing the
# value - otherwise Excel's automatic type-detection would turn them back
# into numbers and silently drop formatting such as leading zeros.
$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r,1).Value = $row[0]

    $q4.Cells.Item($r,2).NumberFormat = "@"
    $q4.Cells.Item($r,2).Value = $row[1]

    $q4.Cells.Item($r,3).Value = $row[2]

    $q4.Cells.Item($r,4).NumberFormat = "@"
    $q4.Cells.Item($r,4).Value = $row[3]

    $q4.Cells.Item($r,5).NumberFormat = "@"
    $q4.Cells.Item($r,5).Value = $row[4]

    $q4.Cells.Item($r,6).NumberFormat = "@"
    $q4.Cells.Item($r,6).Value = $row[5]

    $q4.Cells.Item($r,7).NumberFormat = "@"
    $q4.Cells.Item($r,7).Value = $row[6]

    $q4.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}
